$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2161.6
$ws.Range("I28").Value = 1950.5
$ws.Range("M28").Value = -1465.5
$ws.Range("J28").Value = 3006
$ws.Range("N28").Value = -3976
$ws.Range("K28").Value = 1950.5
$ws.Range("L28").Value = 3006
$ws.Range("I51").Value = 9999
$ws.Range("M51").Value = -9515
$ws.Range("J51").Value = 10433
$ws.Range("N51").Value = -11401
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 10433
$ws.Range("H51").Value = 10259.4
$ws.Range("I92").Value = 9248.75
$ws.Range("K92").Value = 9248.75
$ws.Range("H92").Value = 9248.75
$ws.Range("M92").Value = -8000.75
$ws.Range("H106").Value = 1836.6666
$ws.Range("I106").Value = 1836.6666
$ws.Range("M106").Value = -1205.6666
$ws.Range("K106").Value = 1836.6666
$ws.Range("H125").Value = 5335.4443
$ws.Range("I125").Value = 6977
$ws.Range("M125").Value = -60333
$ws.Range("J125").Value = 4514.6665
$ws.Range("N125").Value = -45551.9985
$ws.Range("K125").Value = 62793
$ws.Range("L125").Value = 40631.9985
$ws.Range("I132").Value = 2108
$ws.Range("M132").Value = -3794
$ws.Range("K132").Value = 6324
$ws.Range("H132").Value = 2108
$ws.Range("K138").Value = 8654.143199999999
$ws.Range("H138").Value = 3142.3215
$ws.Range("M138").Value = -3514.143199999999
$ws.Range("I138").Value = 2884.7144
$ws.Range("J138").Value = 3228.1904
$ws.Range("N138").Value = -19964.5712
$ws.Range("L138").Value = 9684.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("L121").Value = 0
$ws.Range("I132").Value = 1467.3334
$ws.Range("M132").Value = -1872.0002
$ws.Range("K132").Value = 4402.0002
$ws.Range("H132").Value = 1467.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K107").Value = 255.5
$ws.Range("H107").Value = 255.5
$ws.Range("I107").Value = 255.5
$ws.Range("M107").Value = 1664.5
$ws.Range("I134").Value = 1307.2307
$ws.Range("K134").Value = 3921.6921
$ws.Range("H134").Value = 1307.2307
$ws.Range("M134").Value = -1386.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("M16").Value = -1213
$ws.Range("K16").Value = 1500
$ws.Range("H41").Value = 14999.5
$ws.Range("M41").Value = -14572
$ws.Range("I41").Value = 15000
$ws.Range("N41").Value = -15855
$ws.Range("J41").Value = 14999
$ws.Range("L41").Value = 14999
$ws.Range("K41").Value = 15000
$ws.Range("H58").Value = 2013.2858
$ws.Range("I58").Value = 2298.75
$ws.Range("M58").Value = -2095.75
$ws.Range("K58").Value = 2298.75
$ws.Range("L59").Value = 35000
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("N59").Value = -37290
$ws.Range("J68").Value = 60000
$ws.Range("N68").Value = -61498
$ws.Range("L68").Value = 60000
$ws.Range("H68").Value = 60000
$ws.Range("H71").Value = 60000
$ws.Range("N71").Value = -187488
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("M86").Value = -18148069
$ws.Range("J86").Value = 652399
$ws.Range("N86").Value = -654645
$ws.Range("L86").Value = 652399
$ws.Range("K86").Value = 18149192
$ws.Range("H86").Value = 12316928
$ws.Range("I86").Value = 18149192
$ws.Range("K89").Value = 90745960
$ws.Range("L89").Value = 3261995
$ws.Range("H89").Value = 12316928
$ws.Range("I89").Value = 18149192
$ws.Range("M89").Value = -90740344
$ws.Range("J89").Value = 652399
$ws.Range("N89").Value = -3273227
$ws.Range("H99").Value = 5287.625
$ws.Range("M99").Value = -4116.4287
$ws.Range("I99").Value = 5614.4287
$ws.Range("K99").Value = 5614.4287
$ws.Range("N107").Value = -4357.3333
$ws.Range("K107").Value = 314.16666
$ws.Range("L107").Value = 517.3333
$ws.Range("H107").Value = 354.8
$ws.Range("I107").Value = 314.16666
$ws.Range("M107").Value = 1605.83334
$ws.Range("J107").Value = 517.3333
$ws.Range("I113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("K113").Value = 1500
$ws.Range("H113").Value = 1500
$ws.Range("I126").Value = 5614.4287
$ws.Range("M126").Value = -14373.2861
$ws.Range("K126").Value = 16843.2861
$ws.Range("H126").Value = 5287.625
$ws.Range("H136").Value = 2013.2858
$ws.Range("M136").Value = -4346.25
$ws.Range("I136").Value = 2298.75
$ws.Range("K136").Value = 6896.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K3").Value = 1498.5
$ws.Range("H3").Value = 499.5
$ws.Range("I3").Value = 499.5
$ws.Range("M3").Value = -1386.5
$ws.Range("K5").Value = 2307.75
$ws.Range("H5").Value = 815.2
$ws.Range("I5").Value = 769.25
$ws.Range("M5").Value = -2195.75
$ws.Range("K135").Value = 6923.25
$ws.Range("H135").Value = 815.2
$ws.Range("M135").Value = -4388.25
$ws.Range("I135").Value = 769.25
$ws.Range("M140").Value = 2420.9999
$ws.Range("K140").Value = 2759.0001
$ws.Range("H140").Value = 1189.75
$ws.Range("I140").Value = 919.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J33").Value = 19000
$ws.Range("N33").Value = -19504
$ws.Range("L33").Value = 19000
$ws.Range("H33").Value = 15285.714
$ws.Range("H102").Value = 2269
$ws.Range("I102").Value = 2145.9
$ws.Range("M102").Value = -523.9000000000001
$ws.Range("K102").Value = 2145.9
$ws.Range("I126").Value = 6224.75
$ws.Range("N126").Value = -22940
$ws.Range("M126").Value = -16204.25
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 18674.25
$ws.Range("L126").Value = 18000
$ws.Range("H126").Value = 6149.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M46").Value = -4367.5557
$ws.Range("J46").Value = 5555.5557
$ws.Range("N46").Value = -5931.5557
$ws.Range("K46").Value = 4555.5557
$ws.Range("L46").Value = 5555.5557
$ws.Range("H46").Value = 5055.5557
$ws.Range("I46").Value = 4555.5557
$ws.Range("H58").Value = 8131.6665
$ws.Range("I58").Value = 3197
$ws.Range("M58").Value = -2937
$ws.Range("J58").Value = 10599
$ws.Range("N58").Value = -11119
$ws.Range("L58").Value = 10599
$ws.Range("K58").Value = 3197
$ws.Range("H136").Value = 1769.8
$ws.Range("M136").Value = -2586.75
$ws.Range("I136").Value = 1712.25
$ws.Range("N136").Value = -11100
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("K136").Value = 5136.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5648.5
$ws.Range("I62").Value = 6349.25
$ws.Range("M62").Value = -5725.25
$ws.Range("J62").Value = 4247
$ws.Range("N62").Value = -5495
$ws.Range("K62").Value = 6349.25
$ws.Range("L62").Value = 4247
$ws.Range("H65").Value = 5648.5
$ws.Range("I65").Value = 6349.25
$ws.Range("M65").Value = -28626.25
$ws.Range("J65").Value = 4247
$ws.Range("N65").Value = -27475
$ws.Range("K65").Value = 31746.25
$ws.Range("L65").Value = 21235
$ws.Range("I81").Value = 4981.8
$ws.Range("M81").Value = -8902.6
$ws.Range("J81").Value = 2500
$ws.Range("N81").Value = -7122
$ws.Range("K81").Value = 9963.6
$ws.Range("L81").Value = 5000
$ws.Range("H81").Value = 4568.1665
$ws.Range("H84").Value = 4568.1665
$ws.Range("M84").Value = -44514
$ws.Range("J84").Value = 2500
$ws.Range("I84").Value = 4981.8
$ws.Range("N84").Value = -187488
$ws.Range("L84").Value = 25000
$ws.Range("K84").Value = 49818
